$d = $word.ActiveDocument

# Fix #517: a table cell that starts with an empty paragraph followed by the
# paragraph that actually carries the cell's content ends up leaving a
# spurious empty paragraph at the beginning of the cell. Remove that leading
# empty paragraph wherever a cell has more than one paragraph and the first
# one is empty (i.e. contains only the end-of-paragraph mark).
foreach ($table in $d.Tables) {
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        for ($c = 1; $c -le $table.Columns.Count; $c++) {
            $cell = $table.Cell($r, $c)
            $cellParas = $cell.Range.Paragraphs
            if ($cellParas.Count -gt 1) {
                $firstPara = $cellParas.Item(1)
                if ($firstPara.Range.Text -eq [char]13) {
                    $firstPara.Range.Delete()
                }
            }
        }
    }
}
